# Update "想去人数" (F column) counts for a handful of rows that are
# duplicated between the "展览" sheet (1st sheet) and the "全部类型" sheet
# (4th sheet), reflecting refreshed scrape output.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value  = 3491
$ws1.Range("F5").Value  = 3491
$ws1.Range("F22").Value = 4866
$ws1.Range("F26").Value = 5982
$ws1.Range("F30").Value = 320
$ws1.Range("F36").Value = 966
$ws1.Range("F40").Value = 848
$ws1.Range("F41").Value = 939

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F8").Value  = 3491
$ws4.Range("F9").Value  = 3491
$ws4.Range("F26").Value = 4866
$ws4.Range("F30").Value = 5982
$ws4.Range("F34").Value = 320
$ws4.Range("F41").Value = 966
$ws4.Range("F45").Value = 848
$ws4.Range("F46").Value = 939
